# Apply Fgf15-Fgfr2 LR-pair data update (Natmi, following Dr Hou advice).
# The sheet grows from 3 data rows (sCs/FAPs/self only, with Sending cluster
# mislabeled "sCs" for all 3 rows) to a full 2x4 sender x target cross of
# FAPs and sCs against ECs, FAPs, M2 and sCs, with refreshed statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf15"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.029478
$ws.Range("H2").Value = 0.088434
$ws.Range("I2").Value = 0.1535387136874709
$ws.Range("J2").Value = 0.1535387136874709
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6105093333333332
$ws.Range("N2").Value = 1.831528
$ws.Range("O2").Value = 0.1519928013857482
$ws.Range("P2").Value = 0.1519928013857482
$ws.Range("Q2").Value = 0.017996594128
$ws.Range("R2").Value = 0.161969347152
$ws.Range("S2").Value = 0.02333677921452302
$ws.Range("T2").Value = 0.02333677921452302

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf15"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.029478
$ws.Range("H3").Value = 0.088434
$ws.Range("I3").Value = 0.1535387136874709
$ws.Range("J3").Value = 0.1535387136874709
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.333134333333334
$ws.Range("N3").Value = 9.999403000000001
$ws.Range("O3").Value = 0.8298192952305696
$ws.Range("P3").Value = 0.8298192952305695
$ws.Range("Q3").Value = 0.09825413387800001
$ws.Range("R3").Value = 0.8842872049020001
$ws.Range("S3").Value = 0.1274093871827453
$ws.Range("T3").Value = 0.1274093871827453

# Row 4: FAPs -> M2
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf15"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.029478
$ws.Range("H4").Value = 0.088434
$ws.Range("I4").Value = 0.1535387136874709
$ws.Range("J4").Value = 0.1535387136874709
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.0004976666666666667
$ws.Range("N4").Value = 0.001493
$ws.Range("O4").Value = 0.0001238994175731532
$ws.Range("P4").Value = 0.0001238994175731531
$ws.Range("Q4").Value = 0.000014670218
$ws.Range("R4").Value = 0.000132031962
$ws.Range("S4").Value = 0.00001902335720080876
$ws.Range("T4").Value = 0.00001902335720080876

# Row 5: FAPs -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf15"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.029478
$ws.Range("H5").Value = 0.088434
$ws.Range("I5").Value = 0.1535387136874709
$ws.Range("J5").Value = 0.1535387136874709
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.07255766666666667
$ws.Range("N5").Value = 0.217673
$ws.Range("O5").Value = 0.01806400396610915
$ws.Range("P5").Value = 0.01806400396610915
$ws.Range("Q5").Value = 0.002138854898
$ws.Range("R5").Value = 0.019249694082
$ws.Range("S5").Value = 0.002773523933001772
$ws.Range("T5").Value = 0.002773523933001772

# Row 6: sCs -> ECs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fgf15"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1625126666666667
$ws.Range("H6").Value = 0.487538
$ws.Range("I6").Value = 0.846461286312529
$ws.Range("J6").Value = 0.846461286312529
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6105093333333332
$ws.Range("N6").Value = 1.831528
$ws.Range("O6").Value = 0.1519928013857482
$ws.Range("P6").Value = 0.1519928013857482
$ws.Range("Q6").Value = 0.09921549978488887
$ws.Range("R6").Value = 0.892939498064
$ws.Range("S6").Value = 0.1286560221712251
$ws.Range("T6").Value = 0.1286560221712251

# Row 7: sCs -> FAPs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf15"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1625126666666667
$ws.Range("H7").Value = 0.487538
$ws.Range("I7").Value = 0.846461286312529
$ws.Range("J7").Value = 0.846461286312529
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.333134333333334
$ws.Range("N7").Value = 9.999403000000001
$ws.Range("O7").Value = 0.8298192952305696
$ws.Range("P7").Value = 0.8298192952305695
$ws.Range("Q7").Value = 0.5416765488682223
$ws.Range("R7").Value = 4.875088939814001
$ws.Range("S7").Value = 0.7024099080478242
$ws.Range("T7").Value = 0.7024099080478242

# Row 8: sCs -> M2
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf15"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.1625126666666667
$ws.Range("H8").Value = 0.487538
$ws.Range("I8").Value = 0.846461286312529
$ws.Range("J8").Value = 0.846461286312529
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.0004976666666666667
$ws.Range("N8").Value = 0.001493
$ws.Range("O8").Value = 0.0001238994175731532
$ws.Range("P8").Value = 0.0001238994175731531
$ws.Range("Q8").Value = 0.00008087713711111111
$ws.Range("R8").Value = 0.000727894234
$ws.Range("S8").Value = 0.0001048760603723444
$ws.Range("T8").Value = 0.0001048760603723444

# Row 9: sCs -> sCs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf15"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.1625126666666667
$ws.Range("H9").Value = 0.487538
$ws.Range("I9").Value = 0.846461286312529
$ws.Range("J9").Value = 0.846461286312529
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.07255766666666667
$ws.Range("N9").Value = 0.217673
$ws.Range("O9").Value = 0.01806400396610915
$ws.Range("P9").Value = 0.01806400396610915
$ws.Range("Q9").Value = 0.01179153989711111
$ws.Range("R9").Value = 0.106123859074
$ws.Range("S9").Value = 0.01529048003310738
$ws.Range("T9").Value = 0.01529048003310738

